$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction9")

$values = @(0,1,2,5,6,7,8,13,14,15,16,18,19,20,21,25,26,27,29,30,31,32,34,35,36,37,38,39,40,42,43,44,46,47,48,49,50,51,52,53,54,55,56,57,59,60,61,62,63,64,65,67,68,70,71,72,73,74,75)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
